# Insert a new data row directly above the existing row 502. This shifts
# the current rows 502:588 down to 503:589 (matching the target dimension
# A1:R589) and leaves room for the brand-new record at row 502.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(502).Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A502").Value = 3
$ws.Range("B502").Value = "Femacal de La Calera"
$ws.Range("C502").Value = "Coquimbo"
$ws.Range("D502").Value = 44694
$ws.Range("E502").Value = 5
$ws.Range("F502").Value = 100114001
$ws.Range("G502").Value = "Papa"
$ws.Range("H502").Value = "Rosara"
$ws.Range("I502").Value = "1a (cosecha)"
$ws.Range("J502").Value = 540
$ws.Range("K502").Value = 6800
$ws.Range("L502").Value = 7000
$ws.Range("M502").Value = 6904
$ws.Range("N502").Value = "$/saco 25 kilos"
$ws.Range("O502").Value = "Provincia de Melipilla"
$ws.Range("P502").Value = 276
$ws.Range("Q502").Value = 25
$ws.Range("R502").Value = "Hortaliza"
